$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 164 (Hortaliza, Terminal
# Hortofrutícola Agro Chillán - Zanahoria), pushing the existing rows 164-200
# down to 165-201.
$ws.Rows(164).Insert()

$ws.Range("A164").Value = 7
$ws.Range("B164").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C164").Value = "Ñuble"
$ws.Range("D164").Value = 44511
$ws.Range("E164").Value = 16
$ws.Range("F164").Value = 100114013
$ws.Range("G164").Value = "Zanahoria"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 100
$ws.Range("K164").Value = 7500
$ws.Range("L164").Value = 8000
$ws.Range("M164").Value = 7750
$ws.Range("N164").Value = "`$/saco 20 kilos"
$ws.Range("O164").Value = "Provincia de Diguillín"
$ws.Range("P164").Value = 388
$ws.Range("Q164").Value = 20
$ws.Range("R164").Value = "Hortaliza"
